$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 3 (M3)")
$ws.Activate()

# Task "Create Space Trader FX project (1)" (row 9) is now complete, now that
# the JavaFX project has been created — reassign it to Stephen and mark it
# Complete, and record 0 remaining hours in the new column E.
$ws.Range("B9").Value = "Stephen"
$ws.Range("C9").Value = "Complete"
$ws.Range("E9").Value = 0

$ws.Range("E9").Select()
